$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 0.9
$ws.Range("D2").Value = 1.01
$ws.Range("E2").Value = 0.4
$ws.Range("F2").Value = 2.52
$ws.Range("G2").Value = 2.31
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

# Row 3
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 0.87
$ws.Range("D3").Value = 1.01
$ws.Range("E3").Value = 0.4
$ws.Range("F3").Value = 2.86
$ws.Range("G3").Value = 2.96
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 0.84
$ws.Range("D4").Value = 1.01
$ws.Range("E4").Value = 0.4
$ws.Range("F4").Value = 2.88
$ws.Range("G4").Value = 2.81
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0
